$d = $word.ActiveDocument

# Locate the paragraph that carries the site-footer copyright notice.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # The two paragraphs immediately preceding the copyright line are an
    # empty "Normal" paragraph and an empty "page-break-before" paragraph;
    # together with the copyright paragraph itself they form the block
    # that must be removed (the commit dropped this generated footer).
    $prev1 = $target.Previous()
    $prev2 = $prev1.Previous()

    $start = $prev2.Range.Start
    $end = $target.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
